# Apply the authored edits to math_table2.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: "는 닮음꼴이다" -> "는 닮음이다"
$ws.Range("C7").Value = "는 닮음이다"

# Row 10: "는 " (trailing space) -> "는" (no trailing space)
$ws.Range("B10").Value = "는"

# Row 16/17 lost their extra (unused) cell style - clear direct formatting
$ws.Range("A16:B17").ClearFormats()

# Row 19 ("dot" / " " / "닷" - duplicate of row 18 "acute") was removed entirely
$ws.Rows.Item(19).Delete()

# Selection now spans the whole table A1:C18
$ws.Range("A1:C18").Select()
